# Daily update at 8 AM UTC
# Append the next day's win-count row (row 86) to the tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A86").Value = 46035
$ws.Range("B86").Value = 195
$ws.Range("C86").Value = 207
$ws.Range("D86").Value = 193

# Match the date formatting used by the rest of column A (e.g. A85).
$ws.Range("A86").NumberFormat = $ws.Range("A85").NumberFormat
